# EBS_Pollock_Inputs.xlsx - "add in f-based spr calcs" prep edits
#
# 1. Controls!B2 (n_years): 100 -> 300
# 2. Controls!B6 (N_1, starting numbers for age 2s): 1 -> 2
# 3. Move the active selection/tab from Recruitment_Mortality!B5 to
#    Controls!B6 (Controls becomes the active sheet/tab again).

$wb = $excel.ActiveWorkbook

$wsControls = $wb.Worksheets.Item("Controls")

# --- value edits ---
$wsControls.Range("B2").Value = 300
$wsControls.Range("B6").Value = 2

# --- selection / active sheet edits ---
# Activate the Controls sheet (was Recruitment_Mortality) and select B6
# (was B3) so it matches the new saved view state.
$wsControls.Activate()
$wsControls.Range("B6").Select()
